# Update FFXIV Leve market-price/profit figures across all class sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect refreshed Universalis
# price data, per the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(20, 8).Value = 800
$ws.Cells.Item(20, 9).Value = 800
$ws.Cells.Item(20, 11).Value = 800
$ws.Cells.Item(20, 13).Value = -570

$ws.Cells.Item(35, 8).Value = 800
$ws.Cells.Item(35, 9).Value = 800
$ws.Cells.Item(35, 11).Value = 800
$ws.Cells.Item(35, 13).Value = -421

$ws.Cells.Item(62, 8).Value = 5568.9165
$ws.Cells.Item(62, 9).Value = 3832.4285
$ws.Cells.Item(62, 11).Value = 3832.4285
$ws.Cells.Item(62, 13).Value = -3208.4285

$ws.Cells.Item(65, 8).Value = 5568.9165
$ws.Cells.Item(65, 9).Value = 3832.4285
$ws.Cells.Item(65, 11).Value = 19162.1425
$ws.Cells.Item(65, 13).Value = -16042.1425

$ws.Cells.Item(86, 8).Value = 5208.1816
$ws.Cells.Item(86, 9).Value = 4558
$ws.Cells.Item(86, 10).Value = 5750
$ws.Cells.Item(86, 11).Value = 4558
$ws.Cells.Item(86, 12).Value = 5750
$ws.Cells.Item(86, 13).Value = -3435
$ws.Cells.Item(86, 14).Value = -7996

$ws.Cells.Item(89, 8).Value = 5208.1816
$ws.Cells.Item(89, 9).Value = 4558
$ws.Cells.Item(89, 10).Value = 5750
$ws.Cells.Item(89, 11).Value = 22790
$ws.Cells.Item(89, 12).Value = 28750
$ws.Cells.Item(89, 13).Value = -17174
$ws.Cells.Item(89, 14).Value = -39982

$ws.Cells.Item(100, 8).Value = 678.53845
$ws.Cells.Item(100, 9).Value = 568.4167
$ws.Cells.Item(100, 11).Value = 568.4167
$ws.Cells.Item(100, 13).Value = -27.41669999999999

$ws.Cells.Item(132, 8).Value = 1547.2258
$ws.Cells.Item(132, 9).Value = 1518.2307
$ws.Cells.Item(132, 11).Value = 4554.6921
$ws.Cells.Item(132, 13).Value = -2024.6921

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 14402.388
$ws.Cells.Item(32, 10).Value = 21841.266
$ws.Cells.Item(32, 12).Value = 21841.266
$ws.Cells.Item(32, 14).Value = -22415.266

$ws.Cells.Item(82, 8).Value = 0
$ws.Cells.Item(82, 9).Value = 0
$ws.Cells.Item(82, 11).Value = 0
$ws.Cells.Item(82, 13).Value = $null

$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 13).Value = $null

$ws.Cells.Item(101, 8).Value = 60000
$ws.Cells.Item(101, 10).Value = 60000
$ws.Cells.Item(101, 12).Value = 60000
$ws.Cells.Item(101, 14).Value = -66490

$ws.Cells.Item(132, 8).Value = 1498.5
$ws.Cells.Item(132, 9).Value = 1454.6316
$ws.Cells.Item(132, 10).Value = 3999
$ws.Cells.Item(132, 11).Value = 4363.8948
$ws.Cells.Item(132, 12).Value = 11997
$ws.Cells.Item(132, 13).Value = -1833.8948
$ws.Cells.Item(132, 14).Value = -17057

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(37, 8).Value = 764.5
$ws.Cells.Item(37, 9).Value = 500
$ws.Cells.Item(37, 11).Value = 500
$ws.Cells.Item(37, 13).Value = -363

$ws.Cells.Item(99, 8).Value = 1814.1052
$ws.Cells.Item(99, 9).Value = 1634.8182
$ws.Cells.Item(99, 11).Value = 1634.8182
$ws.Cells.Item(99, 13).Value = -136.8181999999999

$ws.Cells.Item(103, 8).Value = 0
$ws.Cells.Item(103, 10).Value = 0
$ws.Cells.Item(103, 12).Value = 0
$ws.Cells.Item(103, 14).Value = $null

$ws.Cells.Item(134, 8).Value = 1005.7895
$ws.Cells.Item(134, 9).Value = 700.64703
$ws.Cells.Item(134, 11).Value = 2101.94109
$ws.Cells.Item(134, 13).Value = 433.0589100000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 300
$ws.Cells.Item(4, 9).Value = 0
$ws.Cells.Item(4, 10).Value = 300
$ws.Cells.Item(4, 11).Value = 0
$ws.Cells.Item(4, 12).Value = 300
$ws.Cells.Item(4, 13).Value = $null
$ws.Cells.Item(4, 14).Value = -524

$ws.Cells.Item(22, 8).Value = 271.14285
$ws.Cells.Item(22, 9).Value = 266.5
$ws.Cells.Item(22, 11).Value = 266.5
$ws.Cells.Item(22, 13).Value = 83.5

$ws.Cells.Item(31, 8).Value = 5243.5
$ws.Cells.Item(31, 9).Value = 4033
$ws.Cells.Item(31, 10).Value = 5416.4287
$ws.Cells.Item(31, 11).Value = 4033
$ws.Cells.Item(31, 12).Value = 5416.4287
$ws.Cells.Item(31, 13).Value = -3738
$ws.Cells.Item(31, 14).Value = -6006.4287

$ws.Cells.Item(34, 8).Value = 5243.5
$ws.Cells.Item(34, 9).Value = 4033
$ws.Cells.Item(34, 10).Value = 5416.4287
$ws.Cells.Item(34, 11).Value = 4033
$ws.Cells.Item(34, 12).Value = 5416.4287
$ws.Cells.Item(34, 13).Value = -3831
$ws.Cells.Item(34, 14).Value = -5820.4287

$ws.Cells.Item(35, 8).Value = 1262.5
$ws.Cells.Item(35, 9).Value = 525
$ws.Cells.Item(35, 10).Value = 2000
$ws.Cells.Item(35, 11).Value = 525
$ws.Cells.Item(35, 12).Value = 2000
$ws.Cells.Item(35, 13).Value = -231
$ws.Cells.Item(35, 14).Value = -2588

$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 9).Value = 0
$ws.Cells.Item(52, 11).Value = 0
$ws.Cells.Item(52, 13).Value = $null

$ws.Cells.Item(107, 8).Value = 573.2308
$ws.Cells.Item(107, 9).Value = 534.8570999999999
$ws.Cells.Item(107, 10).Value = 618
$ws.Cells.Item(107, 11).Value = 534.8570999999999
$ws.Cells.Item(107, 12).Value = 618
$ws.Cells.Item(107, 13).Value = 1385.1429
$ws.Cells.Item(107, 14).Value = -4458

$ws.Cells.Item(134, 8).Value = 2908.5806
$ws.Cells.Item(134, 9).Value = 1924.3914
$ws.Cells.Item(134, 10).Value = 5738.125
$ws.Cells.Item(134, 11).Value = 5773.174199999999
$ws.Cells.Item(134, 12).Value = 17214.375
$ws.Cells.Item(134, 13).Value = -3238.174199999999
$ws.Cells.Item(134, 14).Value = -22284.375

$ws.Cells.Item(140, 8).Value = 45000
$ws.Cells.Item(140, 10).Value = 45000
$ws.Cells.Item(140, 12).Value = 45000
$ws.Cells.Item(140, 14).Value = -55360

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 2309.6667
$ws.Cells.Item(34, 9).Value = 1697.4
$ws.Cells.Item(34, 10).Value = 3075
$ws.Cells.Item(34, 11).Value = 5092.200000000001
$ws.Cells.Item(34, 12).Value = 9225
$ws.Cells.Item(34, 13).Value = -5008.200000000001
$ws.Cells.Item(34, 14).Value = -9393

$ws.Cells.Item(54, 8).Value = 2249.5
$ws.Cells.Item(54, 10).Value = 2249.5
$ws.Cells.Item(54, 12).Value = 6748.5
$ws.Cells.Item(54, 14).Value = -7866.5

$ws.Cells.Item(114, 8).Value = 618.36365
$ws.Cells.Item(114, 10).Value = 775
$ws.Cells.Item(114, 12).Value = 2325
$ws.Cells.Item(114, 14).Value = -8833

$ws.Cells.Item(139, 8).Value = 5310.0586
$ws.Cells.Item(139, 9).Value = 4689.4165
$ws.Cells.Item(139, 11).Value = 14068.2495
$ws.Cells.Item(139, 13).Value = -8928.249500000002

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 13).Value = $null

$ws.Cells.Item(37, 8).Value = 0
$ws.Cells.Item(37, 9).Value = 0
$ws.Cells.Item(37, 11).Value = 0
$ws.Cells.Item(37, 13).Value = $null

$ws.Cells.Item(126, 8).Value = 4827.857
$ws.Cells.Item(126, 9).Value = 4397
$ws.Cells.Item(126, 11).Value = 13191
$ws.Cells.Item(126, 13).Value = -10721

$ws.Cells.Item(132, 8).Value = 4155.6875
$ws.Cells.Item(132, 9).Value = 3847.4285
$ws.Cells.Item(132, 10).Value = 6313.5
$ws.Cells.Item(132, 11).Value = 11542.2855
$ws.Cells.Item(132, 12).Value = 18940.5
$ws.Cells.Item(132, 13).Value = -9012.2855
$ws.Cells.Item(132, 14).Value = -24000.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 6244.5
$ws.Cells.Item(16, 9).Value = 6244.5
$ws.Cells.Item(16, 11).Value = 6244.5
$ws.Cells.Item(16, 13).Value = -6074.5

$ws.Cells.Item(93, 8).Value = 926
$ws.Cells.Item(93, 10).Value = 1498.6666
$ws.Cells.Item(93, 12).Value = 1498.6666
$ws.Cells.Item(93, 14).Value = -3994.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2298.5789
$ws.Cells.Item(81, 9).Value = 1980.7646
$ws.Cells.Item(81, 11).Value = 3961.5292
$ws.Cells.Item(81, 13).Value = -2900.5292

$ws.Cells.Item(84, 8).Value = 2298.5789
$ws.Cells.Item(84, 9).Value = 1980.7646
$ws.Cells.Item(84, 11).Value = 19807.646
$ws.Cells.Item(84, 13).Value = -14503.646

$ws.Cells.Item(95, 8).Value = 43789.6
$ws.Cells.Item(95, 10).Value = 43789.6
$ws.Cells.Item(95, 12).Value = 43789.6
$ws.Cells.Item(95, 14).Value = -49281.6

$ws.Cells.Item(122, 8).Value = 1335.9166
$ws.Cells.Item(122, 9).Value = 1303.2727
$ws.Cells.Item(122, 11).Value = 3909.8181
$ws.Cells.Item(122, 13).Value = -1459.8181

$ws.Cells.Item(128, 8).Value = 49999
$ws.Cells.Item(128, 10).Value = 49999
$ws.Cells.Item(128, 12).Value = 49999
$ws.Cells.Item(128, 14).Value = -59959

$ws.Cells.Item(132, 8).Value = 2451.6667
$ws.Cells.Item(132, 9).Value = 2258.125
$ws.Cells.Item(132, 11).Value = 6774.375
$ws.Cells.Item(132, 13).Value = -4244.375

$ws.Cells.Item(136, 8).Value = 39610.555
$ws.Cells.Item(136, 10).Value = 146500.58
$ws.Cells.Item(136, 12).Value = 439501.74
$ws.Cells.Item(136, 14).Value = -444601.74
